# "Fuel" sheet update for both CH and SIN database
# - updated PEN & CO2 in "Fuel" sheet for various fuels CH (data source: KBOB 2009/1:2016)
# - updated PEN & CO2 for natural gas (NG) in "Fuel" sheet for SIN (ecoinvent 3.4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ELECTRICITY sheet: "Swiss consumer energy mix" (row 3) PEN/CO2/reference
# ---------------------------------------------------------------------------
$wsElec = $wb.Worksheets.Item("ELECTRICITY")
$wsElec.Range("E3").Value = 2.52
$wsElec.Range("F3").Value = 0.028
$wsElec.Range("H3").Value = "KBOB 2009/1:2016, ID 45.020 CH-Verbrauchermix, costs in USD-2015"

# ---------------------------------------------------------------------------
# FUELS sheet: update PEN / CO2 / reference for existing fuels, add Biogas
# ---------------------------------------------------------------------------
$wsFuels = $wb.Worksheets.Item("FUELS")

# Natural Gas / NG
$wsFuels.Range("C2").Value = 1.06
$wsFuels.Range("D2").Formula = "=0.228/3.6"
$wsFuels.Range("F2").Value = "KBOB 2009/1:2016, ID 41.002 Erdgas"

# Electricity / GRID
$wsFuels.Range("C3").Value = 2.52
$wsFuels.Range("D3").Formula = "=0.102/3.6"
$wsFuels.Range("F3").Value = "KBOB 2009/1:2016, ID 45.020 CH-Verbrauchermix"

# Oil / OIL
$wsFuels.Range("C5").Value = 1.23
$wsFuels.Range("D5").Formula = "=0.301/3.6"
$wsFuels.Range("F5").Value = "KBOB 2009/1:2016, ID 41.001 Heizöl"

# Coal / COAL
$wsFuels.Range("C6").Formula = "=(1.2+1.45)/2"
$wsFuels.Range("D6").Formula = "=((0.399+0.439)/2)/3.6"
$wsFuels.Range("F6").Value = "KBOB 2009/1:2016, ID 41.004/41.005 (average)"

# Wood / WOOD
$wsFuels.Range("C7").Value = 0.116
$wsFuels.Range("D7").Formula = "=0.027/3.6"
$wsFuels.Range("F7").Value = "KBOB 2009/1:2016, ID 41.006 Stückholz (average)"

# New row 8: Biogas / BIOGAS
$wsFuels.Range("A7").Copy()
$wsFuels.Range("A8").PasteSpecial(-4122)
$wsFuels.Range("B7").Copy()
$wsFuels.Range("B8").PasteSpecial(-4122)
$wsFuels.Range("C7").Copy()
$wsFuels.Range("C8").PasteSpecial(-4122)
$wsFuels.Range("D7").Copy()
$wsFuels.Range("D8").PasteSpecial(-4122)
$wsFuels.Range("E4").Copy()
$wsFuels.Range("E8").PasteSpecial(-4122)
$wsFuels.Range("F7").Copy()
$wsFuels.Range("F8").PasteSpecial(-4122)

$wsFuels.Range("A8").Value = "Biogas"
$wsFuels.Range("B8").Value = "BIOGAS"
$wsFuels.Range("C8").Value = 0.299
$wsFuels.Range("D8").Formula = "=0.13/3.6"
$wsFuels.Range("E8").ClearContents()
$wsFuels.Range("F8").Value = "KBOB 2009/1:2016, ID 41.009 Biogas"

$wsFuels.Activate()
